$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "rNNJp810"
$ws.Range("B2").Value = 231102296
$ws.Range("C2").Value = "qgfjyfj84"
$ws.Range("D2").Value = "xJQ6&%5s"
$ws.Range("F2").Value = "GezVubln"
$ws.Range("G2").Value = "QDYm"
